$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells, reusing the bold/bordered header style
# that's already applied to A1:AC1 (copy format from AC1, then set text).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (82 wins, 80 losses, 0 ties) for every player
# row (2-47).
$ws.Range("AD2:AD47").Value = 82
$ws.Range("AE2:AE47").Value = 80
$ws.Range("AF2:AF47").Value = 0
